# Update hotel reviews data: fill in the previously-blank English_Reviews_num,
# Local_Rank and Total_Reviews_num figures for the hotel_info row.
#
# These values ("3", "341", "4") look numeric, so a plain .Value assignment
# would be stored as a Number. The source data stores them as text (shared
# strings), so we force text entry the classic Excel way - a leading
# apostrophe - and then clear the resulting "quote prefix" number format back
# to Normal so the cells end up as plain text cells with no special style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("hotel_info")

$ws.Cells.Item(2, 7).Value = "'3"
$ws.Cells.Item(2, 7).Style = "Normal"

$ws.Cells.Item(2, 8).Value = "'341"
$ws.Cells.Item(2, 8).Style = "Normal"

$ws.Cells.Item(2, 9).Value = "'4"
$ws.Cells.Item(2, 9).Style = "Normal"
